# 2021-12-03 add the FID metric
# Adds two new metric columns (FID, Paq2Piq) after the existing LPIPS column
# on Sheet1: header cells O1/P1 + zero-filled data cells O2:P7, then moves
# the active selection to Q6 (one cell to the right of the new last column,
# one row up from the former last row) to match the author's final cursor
# position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric headers (row 1)
$ws.Range("O1").Value = "FID"
$ws.Range("P1").Value = "Paq2Piq"

# New metric data, rows 2-7, both columns filled with 0 just like every
# other metric column in the sheet.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 15).Value = 0   # column O
    $ws.Cells.Item($r, 16).Value = 0   # column P
}

# Match the saved selection/active cell from the authored workbook.
$ws.Range("Q6").Select()
